$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 1094:1095, shifting all existing data (rows
# 1094-1143) down by 2 (to 1096-1145).
$ws.Rows("1094:1095").Insert()

# --- Row 1094: new week's "Pintón" quote ---
$ws.Range("A1094").Value = 5
$ws.Range("B1094").Value = 'Macroferia Regional de Talca'
$ws.Range("C1094").Value = 'Maule'
$ws.Range("D1094").Value = 45147
$ws.Range("E1094").Value = 7
$ws.Range("F1094").Value = 'Fruta'
$ws.Range("G1094").Value = 100108
$ws.Range("H1094").Value = 'Tropicales y subtropicales'
$ws.Range("I1094").Value = 100108006
$ws.Range("J1094").Value = 'Plátano'
$ws.Range("K1094").Value = 'Sin especificar'
$ws.Range("L1094").Value = 'Pintón'
$ws.Range("M1094").Value = 800
$ws.Range("N1094").Value = 15000
$ws.Range("O1094").Value = 15000
$ws.Range("P1094").Value = 15000
$ws.Range("Q1094").Value = '$/caja 20 kilos'
$ws.Range("R1094").Value = 'Ecuador'
$ws.Range("S1094").Value = 750
$ws.Range("T1094").Value = 20

# --- Row 1095: new week's "Primera Pintón" quote ---
$ws.Range("A1095").Value = 5
$ws.Range("B1095").Value = 'Macroferia Regional de Talca'
$ws.Range("C1095").Value = 'Maule'
$ws.Range("D1095").Value = 45147
$ws.Range("E1095").Value = 7
$ws.Range("F1095").Value = 'Fruta'
$ws.Range("G1095").Value = 100108
$ws.Range("H1095").Value = 'Tropicales y subtropicales'
$ws.Range("I1095").Value = 100108006
$ws.Range("J1095").Value = 'Plátano'
$ws.Range("K1095").Value = 'Sin especificar'
$ws.Range("L1095").Value = 'Primera Pintón'
$ws.Range("M1095").Value = 1020
$ws.Range("N1095").Value = 16000
$ws.Range("O1095").Value = 17000
$ws.Range("P1095").Value = 16412
$ws.Range("Q1095").Value = '$/caja 20 kilos'
$ws.Range("R1095").Value = 'Ecuador'
$ws.Range("S1095").Value = 821
$ws.Range("T1095").Value = 20
